$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: replace NOUBAIL MOHAMMED -> JEMAA HORMI, amounts 1000 -> 8000 ---
$ws.Range("A2").Value = "JEMAA HORMI"
$ws.Range("B2").Value = "B219321"
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "225400000805987601012173"
$ws.Range("D2").Value = "KHOURIBGA"
$ws.Range("E2").Value = "CA"
$ws.Range("F2").Value = "Direction régionale"
$ws.Range("G2").Value = "001/RRR/AV1"
$ws.Range("H2").Value = "mensuelle"
$ws.Range("I2").Value = 8000
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 8000

# --- Row 3: ZERNAKH ABDELLAH stays, only contract/amounts change ---
$ws.Range("G3").Value = "001/RRR/AV1"
$ws.Range("I3").Value = 8000
$ws.Range("K3").Value = 8000

# --- Row 4: previously the totals row, now becomes a JEMAA HORMI data row (1000) ---
$ws.Range("A4").Value = "JEMAA HORMI"
$ws.Range("B4").Value = "B219321"
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "225400000805987601012173"
$ws.Range("D4").Value = "KHOURIBGA"
$ws.Range("E4").Value = "CA"
$ws.Range("F4").Value = "Direction régionale"
$ws.Range("G4").Value = "001/RRR/AV1"
$ws.Range("H4").Value = "mensuelle"
$ws.Range("I4").Value = 1000
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 1000

# --- Row 5 (new): ZERNAKH ABDELLAH data row (1000) ---
$ws.Range("A5").Value = "ZERNAKH ABDELLAH"
$ws.Range("B5").Value = "IB19558"
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "145101211406073828000084"
$ws.Range("D5").Value = "MARRAKECH BENI MELLAL"
$ws.Range("E5").Value = "BP"
$ws.Range("F5").Value = "Direction régionale"
$ws.Range("G5").Value = "001/RRR/AV1"
$ws.Range("H5").Value = "mensuelle"
$ws.Range("I5").Value = 1000
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 1000

# --- Row 6 (new): totals row (moved from old row 4), now summing to 18000 ---
$ws.Range("A6").Value = " "
$ws.Range("B6").Value = " "
$ws.Range("C6").Value = " "
$ws.Range("D6").Value = " "
$ws.Range("E6").Value = " "
$ws.Range("F6").Value = " "
$ws.Range("G6").Value = " "
$ws.Range("H6").Value = " "
$ws.Range("I6").Value = 18000
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 18000
